$wb = $excel.ActiveWorkbook

# Resolve worksheets by name
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Sheet1 (LP1912) header updates ---
$ws1.Cells.Item(2, 1).Value = "Última actualización: 19:55:23"
$ws1.Cells.Item(3, 1).Value = "Total filas: 353"

$ws1.Cells.Item(137, 1).Value = '11:47:17'
$ws1.Cells.Item(137, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(137, 4).Value = 50
$ws1.Cells.Item(138, 1).Value = '11:52:01'
$ws1.Cells.Item(138, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(138, 4).Value = 45
$ws1.Cells.Item(160, 1).Value = '12:11:52'
$ws1.Cells.Item(160, 3).Value = '14_ABASTO'
$ws1.Cells.Item(160, 4).Value = 81
$ws1.Cells.Item(161, 1).Value = '11:34:59'
$ws1.Cells.Item(161, 3).Value = '215A_EL PATO'
$ws1.Cells.Item(161, 4).Value = 118
$ws1.Cells.Item(259, 3).Value = '215B_EL PATO'
$ws1.Cells.Item(260, 1).Value = '17:48:33'
$ws1.Cells.Item(260, 3).Value = '215_EL PELIGRO'
$ws1.Cells.Item(260, 4).Value = 0
$ws1.Cells.Item(261, 1).Value = '16:44:12'
$ws1.Cells.Item(261, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(261, 4).Value = 64
$ws1.Cells.Item(302, 1).Value = '18:44:57'
$ws1.Cells.Item(302, 3).Value = '14X44_ABASTO'
$ws1.Cells.Item(302, 4).Value = 32
$ws1.Cells.Item(304, 1).Value = '18:12:30'
$ws1.Cells.Item(304, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(304, 4).Value = 64
$ws1.Cells.Item(326, 1).Value = '19:55:23'
$ws1.Cells.Item(326, 2).Value = '20:06'
$ws1.Cells.Item(326, 4).Value = 11
$ws1.Cells.Item(327, 1).Value = '19:48:11'
$ws1.Cells.Item(327, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(327, 4).Value = 19
$ws1.Cells.Item(328, 2).Value = '20:07'
$ws1.Cells.Item(328, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(328, 4).Value = 75
$ws1.Cells.Item(329, 1).Value = '18:52:02'
$ws1.Cells.Item(329, 2).Value = '20:08'
$ws1.Cells.Item(329, 4).Value = 76
$ws1.Cells.Item(330, 2).Value = '20:09'
$ws1.Cells.Item(330, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(330, 4).Value = 85
$ws1.Cells.Item(331, 1).Value = '18:44:57'
$ws1.Cells.Item(331, 3).Value = '14_ABASTO'
$ws1.Cells.Item(331, 4).Value = 88
$ws1.Cells.Item(332, 1).Value = '18:31:25'
$ws1.Cells.Item(332, 2).Value = '20:12'
$ws1.Cells.Item(332, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(332, 4).Value = 101
$ws1.Cells.Item(333, 1).Value = '18:52:02'
$ws1.Cells.Item(333, 2).Value = '20:13'
$ws1.Cells.Item(333, 3).Value = '14_ABASTO'
$ws1.Cells.Item(333, 4).Value = 81
$ws1.Cells.Item(334, 1).Value = '18:44:57'
$ws1.Cells.Item(334, 2).Value = '20:21'
$ws1.Cells.Item(334, 4).Value = 97
$ws1.Cells.Item(335, 1).Value = '18:31:25'
$ws1.Cells.Item(335, 2).Value = '20:22'
$ws1.Cells.Item(335, 3).Value = '15_ABASTO'
$ws1.Cells.Item(335, 4).Value = 111
$ws1.Cells.Item(336, 1).Value = '18:44:57'
$ws1.Cells.Item(336, 2).Value = '20:30'
$ws1.Cells.Item(336, 4).Value = 106
$ws1.Cells.Item(337, 1).Value = '18:52:02'
$ws1.Cells.Item(337, 2).Value = '20:31'
$ws1.Cells.Item(337, 3).Value = '10_OLMOS'
$ws1.Cells.Item(337, 4).Value = 99
$ws1.Cells.Item(338, 1).Value = '19:35:56'
$ws1.Cells.Item(338, 2).Value = '20:33'
$ws1.Cells.Item(338, 4).Value = 58
$ws1.Cells.Item(339, 1).Value = '19:48:11'
$ws1.Cells.Item(339, 2).Value = '20:34'
$ws1.Cells.Item(339, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(339, 4).Value = 46
$ws1.Cells.Item(340, 1).Value = '19:35:56'
$ws1.Cells.Item(340, 2).Value = '20:42'
$ws1.Cells.Item(340, 4).Value = 67
$ws1.Cells.Item(341, 1).Value = '19:11:45'
$ws1.Cells.Item(341, 2).Value = '20:43'
$ws1.Cells.Item(341, 4).Value = 92
$ws1.Cells.Item(342, 1).Value = '19:48:11'
$ws1.Cells.Item(342, 2).Value = '20:46'
$ws1.Cells.Item(342, 3).Value = '17_ROMERO'
$ws1.Cells.Item(342, 4).Value = 58
$ws1.Cells.Item(343, 1).Value = '19:11:45'
$ws1.Cells.Item(343, 2).Value = '20:47'
$ws1.Cells.Item(343, 4).Value = 96
$ws1.Cells.Item(344, 1).Value = '18:52:02'
$ws1.Cells.Item(344, 2).Value = '20:48'
$ws1.Cells.Item(344, 3).Value = '215B_EL PATO'
$ws1.Cells.Item(344, 4).Value = 116
$ws1.Cells.Item(345, 2).Value = '20:55'
$ws1.Cells.Item(345, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(345, 4).Value = 67
$ws1.Cells.Item(346, 1).Value = '19:55:23'
$ws1.Cells.Item(346, 2).Value = '20:55'
$ws1.Cells.Item(346, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(346, 4).Value = 60
$ws1.Cells.Item(347, 1).Value = '19:48:11'
$ws1.Cells.Item(347, 2).Value = '20:56'
$ws1.Cells.Item(347, 4).Value = 68
$ws1.Cells.Item(348, 2).Value = '20:56'
$ws1.Cells.Item(348, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(348, 4).Value = 105
$ws1.Cells.Item(349, 1).Value = '19:35:56'
$ws1.Cells.Item(349, 2).Value = '20:57'
$ws1.Cells.Item(349, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(349, 4).Value = 82
$ws1.Cells.Item(350, 1).Value = '19:11:45'
$ws1.Cells.Item(350, 2).Value = '21:06'
$ws1.Cells.Item(350, 3).Value = '10_OLMOS'
$ws1.Cells.Item(350, 4).Value = 115
$ws1.Cells.Item(351, 2).Value = '21:07'
$ws1.Cells.Item(351, 3).Value = '10_OLMOS'
$ws1.Cells.Item(351, 4).Value = 79
$ws1.Cells.Item(352, 2).Value = '21:09'
$ws1.Cells.Item(352, 3).Value = '15_ABASTO'
$ws1.Cells.Item(352, 4).Value = 94
$ws1.Cells.Item(353, 1).Value = '19:48:11'
$ws1.Cells.Item(353, 2).Value = '21:10'
$ws1.Cells.Item(353, 3).Value = '15_ABASTO'
$ws1.Cells.Item(353, 4).Value = 82
$ws1.Cells.Item(354, 1).Value = '19:35:56'
$ws1.Cells.Item(354, 2).Value = '21:28'
$ws1.Cells.Item(354, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(354, 4).Value = 113
$ws1.Cells.Item(355, 1).Value = '19:35:56'
$ws1.Cells.Item(355, 2).Value = '21:33'
$ws1.Cells.Item(355, 3).Value = '84_COLONIA URQUIZA-ESC 49'

# --- New rows appended at the bottom of sheet1 ---
$ws1.Cells.Item(356, 1).Value = "19:48:11"
$ws1.Cells.Item(356, 2).Value = "21:34"
$ws1.Cells.Item(356, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(356, 4).Value = 106
$ws1.Cells.Item(356, 5).Value = "LP1912"

$ws1.Cells.Item(357, 1).Value = "19:55:23"
$ws1.Cells.Item(357, 2).Value = "21:34"
$ws1.Cells.Item(357, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(357, 4).Value = 99
$ws1.Cells.Item(357, 5).Value = "LP1912"

$ws1.Cells.Item(358, 1).Value = "19:48:11"
$ws1.Cells.Item(358, 2).Value = "21:46"
$ws1.Cells.Item(358, 3).Value = "14X44_ABASTO"
$ws1.Cells.Item(358, 4).Value = 118
$ws1.Cells.Item(358, 5).Value = "LP1912"

# --- Sheet2 (LP1912-215) and Sheet3 (6203-6173) header updates ---
$ws2.Cells.Item(2, 1).Value = "Última actualización: 19:55:23"
$ws3.Cells.Item(2, 1).Value = "Última actualización: 19:55:23"
